$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "rechazada"/"38.46%" -> "pendiente"/"100.00%"
# Force the percentage-looking text to stay as text (matching the original
# file, where these percentages are stored as literal strings, not numbers)
# instead of Excel's automatic number/percent conversion.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("A2").Value = "pendiente"
$ws.Range("B2").Value = "100.00%"
$ws.Range("B2").ClearFormats()

# Row 3 ("aceptada"/"61.54%") is removed entirely, shrinking the table
# from A1:B3 down to A1:B2.
$ws.Rows(3).Delete()
